$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- RUNMANAGER (sheet1): add a new test case row ---
$ws1.Range("A6").Value = "postProductByReadingRequestFromFile"
$ws1.Range("B6").Value = "Post a new product by reading from a json file"
$ws1.Range("C6").Value = "Yes"
$ws1.Range("D6").Value = "'1"
$ws1.Range("E6").Value = "'1"

# Existing "No" -> "Yes" for the two getCountryDetailsTest / getProducts rows
$ws1.Range("C2").Value = "Yes"
$ws1.Range("C3").Value = "Yes"

# --- TESTDATA (sheet2): flip Execute flag for rows 8 & 10, add new row 11 ---
$ws2.Range("B8").Value = "Yes"
$ws2.Range("B10").Value = "Yes"

$ws2.Range("A11").Value = "postProductByReadingRequestFromFile"
$ws2.Range("B11").Value = "Yes"
$ws2.Range("C11").Value = "'"
$ws2.Range("D11").Value = "'"
$ws2.Range("E11").Value = "'"
$ws2.Range("F11").Value = "'checking"
$ws2.Range("G11").Value = "'"

# widen column A on TESTDATA to fit the new (longer) test case name
# (ColumnWidth="34" renders as width=34 chars once written to XML; Excel's
# COM API adds ~5/6 of a character internally so we back that out here)
$ws2.Columns.Item(1).ColumnWidth = 33.1666666666667

# --- selection / active sheet bookkeeping ---
# Final state: RUNMANAGER tab selected with C3 active; TESTDATA not selected,
# with B4 as its remembered active cell.
$ws2.Range("B4").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C3").Select() | Out-Null
